# Updates the cryptos list (prices and 1h volume changes) to the latest
# values pulled on Fri Mar  1 13:24:01 UTC 2024, and fixes the ordering of
# the Monero / ARBITRUM rows (41 and 42) which had swapped coin/link data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    # Forces the written value to be stored as text, even when it looks
    # like a number (e.g. "407.40"), matching the original sheet where the
    # whole Price column is plain text.
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "62.546.83"
$ws.Range("E2").Value = "  -0.29%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.440.97"
$ws.Range("E3").Value = "  -0.83%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.13%  "

# Row 5 - BNB
Set-TextCell "D5" "407.40"
$ws.Range("E5").Value = "  -0.30%  "

# Row 6 - Solana
Set-TextCell "D6" "134.07"
$ws.Range("E6").Value = "  +2.38%  "

# Row 7 - XRP
Set-TextCell "D7" "0.593"
$ws.Range("E7").Value = "  -0.82%  "

# Row 8 - USDC
Set-TextCell "D8" "0.999"
$ws.Range("E8").Value = "  -0.20%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.686"
$ws.Range("E9").Value = "  -0.61%  "

# Row 10 - Dogecoin
Set-TextCell "D10" "0.123"
$ws.Range("E10").Value = "  -3.58%  "

# Row 11 - Avalanche
Set-TextCell "D11" "42.42"
$ws.Range("E11").Value = "  -0.64%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.61%  "

# Row 13 - Polkadot
Set-TextCell "D13" "8.46"
$ws.Range("E13").Value = "  -2.88%  "

# Row 14 - Chainlink
$ws.Range("E14").Value = "  -0.64%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.449.43"
$ws.Range("E15").Value = "  -0.02%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "62.428.34"
$ws.Range("E16").Value = "  -0.51%  "

# Row 17 - Uniswap
$ws.Range("E17").Value = "  +5.21%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  -2.21%  "

# Row 19 - ShibaInu
Set-TextCell "D19" "0.0000132"
$ws.Range("E19").Value = "  -2.70%  "

# Row 20 - ImmutableX
$ws.Range("E20").Value = "  -4.94%  "

# Row 21 - Litecoin
Set-TextCell "D21" "84.23"
$ws.Range("E21").Value = "  +2.02%  "

# Row 22 - BitcoinCash
Set-TextCell "D22" "315.09"
$ws.Range("E22").Value = "  +1.99%  "

# Row 23 - InternetComputer(DFINITY)
Set-TextCell "D23" "12.96"
$ws.Range("E23").Value = "  -1.28%  "

# Row 24 - PancakeSwap
Set-TextCell "D24" "3.17"
$ws.Range("E24").Value = "  -0.53%  "

# Row 25 - LEO
Set-TextCell "D25" "4.72"
$ws.Range("E25").Value = "  +7.92%  "

# Row 26 - EthereumClassic
Set-TextCell "D26" "29.82"
$ws.Range("E26").Value = "  -1.64%  "

# Row 27 - Filecoin
Set-TextCell "D27" "8.26"
$ws.Range("E27").Value = "  -0.17%  "

# Row 28 - RenderToken
Set-TextCell "D28" "7.61"
$ws.Range("E28").Value = "  -1.90%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  +2.42%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  -3.95%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -3.49%  "

# Row 32 - Dai (unchanged)

# Row 33 - InjectiveProtocol
Set-TextCell "D33" "42.26"
$ws.Range("E33").Value = "  -2.41%  "

# Row 34 - Cosmos
Set-TextCell "D34" "11.40"
$ws.Range("E34").Value = "  -4.40%  "

# Row 35 - VeChain
$ws.Range("E35").Value = "  -1.56%  "

# Row 36 - OKB
Set-TextCell "D36" "51.60"
$ws.Range("E36").Value = "  -1.81%  "

# Row 37 - FirstDigitalUSD
Set-TextCell "D37" "0.998"
$ws.Range("E37").Value = "  -0.03%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  -4.22%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  -1.74%  "

# Row 40 - TheGraph
$ws.Range("E40").Value = "  +11.13%  "

# Row 41 / 42 - Monero and ARBITRUM had their Coin/Link columns swapped;
# row 41 now holds ARBITRUM's data and row 42 holds Monero's data.
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D41" "1.99"
$ws.Range("E41").Value = "  +0.10%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D42" "137.59"
# E42 stays "  -0.19%  " (unchanged)

# Row 43 - Stellar
$ws.Range("E43").Value = "  -0.46%  "

# Row 44 - NEARProtocol
$ws.Range("E44").Value = "  +1.71%  "

# Row 45 - Celestia
Set-TextCell "D45" "16.86"
$ws.Range("E45").Value = "  -3.85%  "

# Row 46 - WEMIXToken
$ws.Range("E46").Value = "  -1.11%  "

# Row 47 - EnergySwap
Set-TextCell "D47" "21.34"
$ws.Range("E47").Value = "  -4.69%  "

# Row 48 - Maker
$ws.Range("D48").Value = "2.129.68"
$ws.Range("E48").Value = "  -3.49%  "

# Row 49 - ApeXProtocol
$ws.Range("E49").Value = "  -4.01%  "

# Row 50 - ThetaToken
$ws.Range("E50").Value = "  +2.97%  "

# Row 51 - Fetch.AI
$ws.Range("E51").Value = "  +22.08%  "
